$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# "EmployeeID" -> "Employee ID"
$ws.Range("B1").Value = "Employee ID"

# --- Replace the employee rows (old rows 2-12) with the new, shorter
#     table (new rows 2-4). Values are written as text so they keep the
#     existing "text" cell style (s="1", numFmtId 49) and are stored as
#     shared strings, matching the target workbook.
$ws.Range("A2").Value = "Louise"
$ws.Range("B2").Value = "124152"

$ws.Range("A3").Value = "John Marston"
$ws.Range("B3").Value = "251241"

$ws.Range("A4").Value = "Dmitry"
$ws.Range("B4").Value = "123241"

# Remove the old trailing rows (5-12) so the used range / dimension
# shrinks back down to A1:B4.
$ws.Range("A5:B12").Clear()

# --- Column B width -----------------------------------------------------
# Widen column B to fit the longest new entry ("John Marston") - mirrors
# the author's bestFit/AutoFit of column B after editing the data.
# (ColumnWidth is snapped to this host's pixel grid; 21.6 lands on the
# closest attainable width to Excel's computed best-fit of 22.43 chars.)
$ws.Columns.Item(2).EntireColumn.AutoFit()
$ws.Columns.Item(2).ColumnWidth = 21.6

# --- Selection state ------------------------------------------------------
# Leftover UI selection landed on B9 when the author saved.
$ws.Range("B9").Select()
